# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The workbook is a "Estado de Cuenta" (account statement) for NIT 9017183432.
# This edit trims the detail table down to a single worker / single period
# (WILFRIDO PUELLO GARCIA, period 2409) by removing the rows that belonged to
# the second worker (LUIS EDUARDO AUSAQUE RODRIGUEZ, periods 2502/2501/2412/
# 2411/2410), updates the summary totals accordingly, and lets the trailing
# signature rows collapse upward into the freed space.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary block (rows 7-13) ------------------------------------------
# "VALOR MORA" total for the NIT drops to just WILFRIDO's outstanding value.
$ws.Range("E11").Value = 64000

# "Cant. Trabajadores" (worker count) and "Cant. Periodos" (period count)
# now both reflect the single remaining worker / single remaining period.
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1

# --- Detail table (rows 15-21) ------------------------------------------
# Row 16 (WILFRIDO PUELLO GARCIA / CC 3829906 / period 2409) stays as-is.
# Rows 17-21 held the second worker's five overdue periods; delete those
# rows entirely so everything below (the signature/footer block) shifts up.
$ws.Rows("17:21").Delete()

# Column D ("Nombre Trabajador") no longer needs to fit the longer name
# ("LUIS EDUARDO AUSAQUE RODRIGUEZ"), so its best-fit width shrinks back
# down to whatever the remaining names require.
$ws.Columns("D:D").AutoFit()

Write-Output "done"
